# SYMB Material List.xlsx - add "Include In Material List" column to Table1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the table to include the new 11th column (K) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K3"))

# Set the new header cell content (also registers the new shared string)
$ws.Range("K1").Value = "Include In Material List"

# Give the new table column its real name (overwrites the auto "Column11" default)
$tbl.ListColumns.Item(11).Name = "Include In Material List"

# --- Column width touch-ups to mirror the resulting autofit pass ---
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 21.666666666666668
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668
$ws.Columns.Item(5).ColumnWidth = 19.5
$ws.Columns.Item(6).ColumnWidth = 17.333333333333332
$ws.Columns.Item(7).ColumnWidth = 22.333333333333332
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666
$ws.Columns.Item(10).ColumnWidth = 12.666666666666666
$ws.Columns.Item(11).ColumnWidth = 30.166666666666668

# --- Move the active selection to K2, matching the saved view state ---
$ws.Range("K2").Select() | Out-Null
